# Apply timesheet updates to the "Zeitplanung" sheet, in the "Implementation"
# block (rows 19-21: "Vorbereitungen ...", "Schulung, Wissen aneignen" and
# "Anforderung #01"), adding logged hours for a few more days/columns.
# Dependent formulas (row/column totals, grand total, and the derived
# "Ist Arbeitszeit - Übersicht" sheet) recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Zeitplanung")

# Row 19
$ws.Range("W19").Value = 6
$ws.Range("X19").Value = 7
$ws.Range("Y19").Value = 7
$ws.Range("AH19").Value = 4
$ws.Range("AR19").Value = 7
$ws.Range("AS19").Value = 5
$ws.Range("AU19").Value = 0.5

# Row 20
$ws.Range("W20").Value = 2
$ws.Range("X20").Value = 1
$ws.Range("Y20").Value = 1
$ws.Range("AR20").Value = 1
$ws.Range("AT20").Value = 3

# Row 21
$ws.Range("AS21").Value = 3
$ws.Range("AT21").Value = 5
$ws.Range("AU21").Value = 0.5

# Recalculate so dependent formulas (row/column sums, grand totals and the
# linked "Ist Arbeitszeit - Übersicht" sheet) update.
$excel.CalculateFullRebuild()

# Restore the view state: scroll position and active selection.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$win.ScrollColumn = 1
$ws.Range("AJ49").Select()
